# Add missing ISO's to RCP mapping file
# New countries/territories appended as rows 236-244 on "EDGAR32 & IEA" (sheet3):
#   cuw/Curacao, mne/Montenegro, pse/Palestine, rou/Romania, srb/Serbia,
#   srb (kosovo)/Kosovo, ssd/South Sudan, sxm/Other S. & Cent. America, tls/Timor-Leste

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDGAR32 & IEA")

# --- Column A (iso codes) first, top to bottom, so new shared strings are
#     appended in the same order the original author typed them. ---
$ws.Range("A236").Value = "cuw"
$ws.Range("A237").Value = "mne"
$ws.Range("A238").Value = "pse"
$ws.Range("A239").Value = "rou"
$ws.Range("A240").Value = "srb"
$ws.Range("A241").Value = "srb (kosovo)"
$ws.Range("A242").Value = "ssd"
$ws.Range("A243").Value = "sxm"
$ws.Range("A244").Value = "tls"

# --- Column B (country / region names) next, top to bottom. ---
$ws.Range("B236").Value = "Curacao"
$ws.Range("B237").Value = "Montenegro"
$ws.Range("B238").Value = "Palestine"
$ws.Range("B239").Value = "Romania"
$ws.Range("B240").Value = "Serbia"
$ws.Range("B241").Value = "Kosovo"
$ws.Range("B242").Value = "South Sudan"
$ws.Range("B243").Value = "Other S. & Cent. America"
$ws.Range("B244").Value = "Timor-Leste"

# --- Column C (EDGAR32 lookup code) - blank for row 238 (Palestine). ---
$ws.Range("C236").Value = "OTHERLATIN"
$ws.Range("C237").Value = "SERBMONT"
$ws.Range("C239").Value = "ROMANIA"
$ws.Range("C240").Value = "SERBMONT"
$ws.Range("C241").Value = "SERBMONT"
$ws.Range("C242").Value = "SUDAN"
$ws.Range("C243").Value = "OTHERLATIN"
$ws.Range("C244").Value = "OTHERASIA"

# --- Column D (RCP Template Reg #) ---
$ws.Range("D236").Value = 24
$ws.Range("D237").Value = 33
$ws.Range("D238").Value = 18
$ws.Range("D239").Value = 33
$ws.Range("D240").Value = 33
$ws.Range("D241").Value = 33
$ws.Range("D242").Value = 21
$ws.Range("D243").Value = 24
$ws.Range("D244").Value = 32

# --- Column F (Country lookup) - mirrors column C, except Palestine -> Jordan ---
$ws.Range("F236").Value = "OTHERLATIN"
$ws.Range("F237").Value = "SERBMONT"
$ws.Range("F238").Value = "JORDAN"
$ws.Range("F239").Value = "ROMANIA"
$ws.Range("F240").Value = "SERBMONT"
$ws.Range("F241").Value = "SERBMONT"
$ws.Range("F242").Value = "SUDAN"
$ws.Range("F243").Value = "OTHERLATIN"
$ws.Range("F244").Value = "OTHERASIA"

# --- Column H (IMAGE24 Reg) ---
$ws.Range("H236").Value = "Rest Central America"
$ws.Range("H237").Value = "Central Europe"
$ws.Range("H238").Value = "Middle East"
$ws.Range("H239").Value = "Eastern Africa"
$ws.Range("H240").Value = "Central Europe"
$ws.Range("H241").Value = "Central Europe"
$ws.Range("H242").Value = "Eastern Africa"
$ws.Range("H243").Value = "Rest Central America"
$ws.Range("H244").Value = "Oceania"

# --- Carry over the shaded/striped formatting used by the other
#     "OTHERLATIN"/"OTHERASIA" rows onto the matching new cells. ---
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C236").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D236").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F236").PasteSpecial(-4122) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("C243").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D243").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F243").PasteSpecial(-4122) | Out-Null

$ws.Range("C142").Copy() | Out-Null
$ws.Range("C244").PasteSpecial(-4122) | Out-Null
$ws.Range("D142").Copy() | Out-Null
$ws.Range("D244").PasteSpecial(-4122) | Out-Null
$ws.Range("F142").Copy() | Out-Null
$ws.Range("F244").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- This sheet is the one the author ended up looking at / saving from:
#     make it the active tab and leave the selection on the last new row. ---
$ws.Activate()
$ws.Range("C247").Select()
